$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet is a weekly log of price observations for "Acelga" (chard) at
# Vega Central Mapocho de Santiago. Each week contributes a block of 3 rows
# (one per quality grade: Extra / Primera / Segunda). A new week's block
# (fecha serial 44476, 2021-10-07) is being added, inserted right above the
# block that used to start at row 288, pushing that row and everything below
# it down by 3 rows (old 288:335 -> new 291:338).

# Make room for the new week block (inherits formatting from the row above,
# i.e. the same date-formatted style used throughout column D).
$ws.Rows("288:290").Insert()

# Row 288: Acelga / Extra, fecha 44476 (2021-10-07)
$ws.Range("A288").Value = 9
$ws.Range("B288").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C288").Value = "Metropolitana"
$ws.Range("D288").Value = 44476
$ws.Range("E288").Value = 13
$ws.Range("F288").Value = 100112009
$ws.Range("G288").Value = "Acelga"
$ws.Range("H288").Value = "Sin especificar"
$ws.Range("I288").Value = "Extra"
$ws.Range("J288").Value = 43
$ws.Range("K288").Value = 12000
$ws.Range("L288").Value = 13000
$ws.Range("M288").Value = 12488
$ws.Range("N288").Value = "`$/docena de atados"
$ws.Range("O288").Value = "Región Metropolitana"
$ws.Range("P288").Value = 4163
$ws.Range("Q288").Value = 3
$ws.Range("R288").Value = "Hortaliza"

# Row 289: Acelga / Primera, fecha 44476 (2021-10-07)
$ws.Range("A289").Value = 9
$ws.Range("B289").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C289").Value = "Metropolitana"
$ws.Range("D289").Value = 44476
$ws.Range("E289").Value = 13
$ws.Range("F289").Value = 100112009
$ws.Range("G289").Value = "Acelga"
$ws.Range("H289").Value = "Sin especificar"
$ws.Range("I289").Value = "Primera"
$ws.Range("J289").Value = 61
$ws.Range("K289").Value = 10000
$ws.Range("L289").Value = 11000
$ws.Range("M289").Value = 10508
$ws.Range("N289").Value = "`$/docena de atados"
$ws.Range("O289").Value = "Región Metropolitana"
$ws.Range("P289").Value = 3503
$ws.Range("Q289").Value = 3
$ws.Range("R289").Value = "Hortaliza"

# Row 290: Acelga / Segunda, fecha 44476 (2021-10-07)
$ws.Range("A290").Value = 9
$ws.Range("B290").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C290").Value = "Metropolitana"
$ws.Range("D290").Value = 44476
$ws.Range("E290").Value = 13
$ws.Range("F290").Value = 100112009
$ws.Range("G290").Value = "Acelga"
$ws.Range("H290").Value = "Sin especificar"
$ws.Range("I290").Value = "Segunda"
$ws.Range("J290").Value = 34
$ws.Range("K290").Value = 8000
$ws.Range("L290").Value = 9000
$ws.Range("M290").Value = 8500
$ws.Range("N290").Value = "`$/docena de atados"
$ws.Range("O290").Value = "Región Metropolitana"
$ws.Range("P290").Value = 2833
$ws.Range("Q290").Value = 3
$ws.Range("R290").Value = "Hortaliza"
